$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# --- Fix up urbansim_path column (J) for existing rows: "census_petrale" -> na ---
$ws.Range("J2").Value = "na"
$ws.Range("J3").Value = "na"
$ws.Range("J4").Value = "na"

# --- Row 5: 2035_TM160_NGFr2_NP03_Path1_01 ---
$ws.Range("A5").Value = "NextGenFwys"
$ws.Range("B5").Value = 2035
$ws.Range("C5").Value = "2035_TM160_NGFr2_NP03_Path1_01"
$ws.Range("D5").Value = "NGF_Round2"
$ws.Range("E5").Value = "P1_AllLaneTolling"
$ws.Range("F5").Value = "P1 initial run -  tolls not in 2000$"
$ws.Range("H5").Value = "R2_ALT"
$ws.Range("I5").Value = "https://app.asana.com/0/1203644633064654/1206539107762749/f"
$ws.Range("J5").Value = "na"
$ws.Range("K5").Value = "na"
$ws.Range("K5").Style = "Normal"

# --- Row 6: 2035_TM160_NGFr2_NP04_Path1_01 ---
$ws.Range("A6").Value = "NextGenFwys"
$ws.Range("B6").Value = 2035
$ws.Range("C6").Value = "2035_TM160_NGFr2_NP04_Path1_01"
$ws.Range("D6").Value = "NGF_Round2"
$ws.Range("E6").Value = "P1_AllLaneTolling"
$ws.Range("F6").Value = "P1 initial run -  carpool and tolls.csv fixes"
$ws.Range("H6").Value = "NGF_Networks_NGFround2_P1_01"
$ws.Range("I6").Value = "https://app.asana.com/0/1203644633064654/1207102772074759/f"
$ws.Range("J6").Value = "na"
$ws.Range("K6").Value = "na"
$ws.Range("K6").Style = "Normal"

# --- Row 7: 2035_TM160_NGFr2_NP04_Path1_02 ---
$ws.Range("A7").Value = "NextGenFwys"
$ws.Range("B7").Value = 2035
$ws.Range("C7").Value = "2035_TM160_NGFr2_NP04_Path1_02"
$ws.Range("D7").Value = "NGF_Round2"
$ws.Range("E7").Value = "P1_AllLaneTolling"
$ws.Range("F7").Value = "P1 initial run -  complete street fixes"
$ws.Range("G7").Value = "current"
$ws.Range("H7").Value = "NGF_Networks_NGFround2_P1_01"
$ws.Range("I7").Value = "https://app.asana.com/0/1203644633064654/1207177908266953/f"
$ws.Range("J7").Value = "na"
$ws.Range("K7").Value = "na"
$ws.Range("K7").Style = "Normal"

# --- Row 8: 2035_TM160_NGFr2_NP04_Path1_02_MinTollOff ---
$ws.Range("A8").Value = "NextGenFwys"
$ws.Range("B8").Value = 2035
$ws.Range("C8").Value = "2035_TM160_NGFr2_NP04_Path1_02_MinTollOff"
$ws.Range("D8").Value = "NGF_Round2"
$ws.Range("E8").Value = "P1_AllLaneTolling"
$ws.Range("F8").Value = "P1 initial run -  test with min toll turned off"
$ws.Range("H8").Value = "NGF_Networks_NGFround2_P1_01"
$ws.Range("I8").Value = "https://app.asana.com/0/1201809392759895/1207151709274835/f"
$ws.Range("J8").Value = "na"
$ws.Range("K8").Value = "na"
$ws.Range("K8").Style = "Normal"

# New rows use plain default style except column B which is center-aligned (matches style index 10)
$ws.Range("B5:B8").HorizontalAlignment = -4108

# --- Update view / selection to match the saved workbook state ---
$ws.Range("F18").Select()

# --- Remove the autoFilter (table header filter arrows) ---
$ws.AutoFilterMode = $false
